$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the "Periodo Mora" column (E16:E19) with the new period labels.
$ws.Range("E16").Value = "2303"
$ws.Range("E17").Value = "2304"
$ws.Range("E18").Value = "2305"
$ws.Range("E19").Value = "2306"

# Update "Valor Mora" column (G16:G19) with the new amounts.
$ws.Range("G16").Value = 1533500
$ws.Range("G17").Value = 1533500
$ws.Range("G18").Value = 1533500
$ws.Range("G19").Value = 1533500
